$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in / clear values in column F (and a couple C) for rows 2-25 ---
# (Re-run of imputation / error calc: some previously-filled cells become
# blank, some previously-blank cells get a computed value.)

$ws.Cells.Item(6, 6).Value = 16.43      # F6  (RM 21)  blank -> 16.43
$ws.Cells.Item(8, 6).ClearContents()    # F8  (RM 38)  17.05 -> blank
$ws.Cells.Item(12, 6).Value = 17.45     # F12 (RM 81)  blank -> 17.45
$ws.Cells.Item(14, 6).ClearContents()   # F14 (RM 90)  17.76 -> blank
$ws.Cells.Item(17, 6).Value = 17.78     # F17 (RM 116) blank -> 17.78
$ws.Cells.Item(18, 6).Value = 18.35     # F18 (RM 120) blank -> 18.35
$ws.Cells.Item(19, 6).ClearContents()   # F19 (RM 125) 17.81 -> blank
$ws.Cells.Item(20, 6).ClearContents()   # F20 (RM 134) 17.73 -> blank
$ws.Cells.Item(23, 6).Value = 16.48     # F23 (RM 140) blank -> 16.48

# --- Remove two whole rows: "RM 232" (row 26) and "SC 92" (old row 28,
# which becomes row 27 once RM 232 is removed) ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Remaining cell edits on the rows that shifted up into 26-33 ---
# Row 27 = SC 101
$ws.Cells.Item(27, 3).Value = 10        # C27 blank -> 10
$ws.Cells.Item(27, 6).ClearContents()   # F27 17 -> blank

# Row 28 = SC 105
$ws.Cells.Item(28, 3).ClearContents()   # C28 11.1 -> blank

# Row 29 = SC 119
$ws.Cells.Item(29, 3).ClearContents()   # C29 11.2 -> blank

# Row 30 = SC 120
$ws.Cells.Item(30, 3).Value = 11.4      # C30 blank -> 11.4

# Row 32 = SC 193
$ws.Cells.Item(32, 3).ClearContents()   # C32 10.5 -> blank
